# Add a new log entry (row 47) to the CIS598 log book, mirroring the
# existing Start/Stop/Delta/Notes pattern used by the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start time carries over from the previous entry's Stop time (6:00PM),
# new entry runs to 7:00PM for a 60 minute delta.
$ws.Range("B47").Value = "6:00PM"
$ws.Range("C47").Value = "7:00PM"
$ws.Range("E47").Value = 60
$ws.Range("F47").Value = "Connected new client click action with insert statement into sql db; db now creates new clients and doesn't crash when the unique name constraint is violated. Need to put popup dialog after new client button click to show user the client has succcessfuly been added or not"

# Update the active selection to reflect where the user's cursor ended up
# after entering the new row of data.
$ws.Range("F47").Select()
